# Append: 2025-10-25 18:22 JST
# Update the "取得日時" (acquired datetime) timestamp in column A for the
# existing data rows (2-13) on the "ランサーズ" sheet from
# "2025-10-25 12:42:53" to "2025-10-25 18:22:19".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-25 18:22:19"

for ($row = 2; $row -le 13; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
